$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: curated dimensions -> measures (sector-de-actividad, sexo are no
# longer dimensions in the refreshed metadata)
$ws.Range("A2").Value = "iaest-measure:sector-de-actividad"
$ws.Range("F2").Value = "iaest-measure:sexo"

# Row 3: role changes from "dim" to "medida" for those same two columns
$ws.Range("A3").Value = "medida"
$ws.Range("F3").Value = "medida"

# Row 4: datatype changes from "skos:Concept" to "xsd:int" for those columns
$ws.Range("A4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"

# Row 5 held the (now obsolete) mapping-file references; drop the whole row
$ws.Rows.Item(5).Delete()
